$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.055887666666667
$ws.Range("H2").Value = 12.167663
$ws.Range("I2").Value = 0.4763357569530485
$ws.Range("J2").Value = 0.4763357569530485
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.9419
$ws.Range("N2").Value = 44.8257
$ws.Range("O2").Value = 0.9802815267721257
$ws.Range("P2").Value = 0.9802815267721257
$ws.Range("Q2").Value = 60.60266792656667
$ws.Range("R2").Value = 545.4240113391
$ws.Range("S2").Value = 0.4669431430820906
$ws.Range("T2").Value = 0.4669431430820906
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.055887666666667
$ws.Range("H3").Value = 12.167663
$ws.Range("I3").Value = 0.4763357569530485
$ws.Range("J3").Value = 0.4763357569530485
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.141981
$ws.Range("N3").Value = 0.425943
$ws.Range("O3").Value = 0.009314836229169864
$ws.Range("P3").Value = 0.009314836229169864
$ws.Range("Q3").Value = 0.575858986801
$ws.Range("R3").Value = 5.182730881209
$ws.Range("S3").Value = 0.004436989566115307
$ws.Range("T3").Value = 0.004436989566115307
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.055887666666667
$ws.Range("H4").Value = 12.167663
$ws.Range("I4").Value = 0.4763357569530485
$ws.Range("J4").Value = 0.4763357569530485
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.158577
$ws.Range("N4").Value = 0.475731
$ws.Range("O4").Value = 0.01040363699870454
$ws.Range("P4").Value = 0.01040363699870454
$ws.Range("Q4").Value = 0.643170498517
$ws.Range("R4").Value = 5.788534486653001
$ws.Range("S4").Value = 0.00495562430484267
$ws.Range("T4").Value = 0.004955624304842669
$ws.Range("I5").Value = 0.4564376967244237
$ws.Range("J5").Value = 0.4564376967244237
$ws.Range("M5").Value = 14.9419
$ws.Range("N5").Value = 44.8257
$ws.Range("O5").Value = 0.9802815267721257
$ws.Range("P5").Value = 0.9802815267721257
$ws.Range("Q5").Value = 58.07110165463333
$ws.Range("R5").Value = 522.6399148917
$ws.Range("S5").Value = 0.4474374422213705
$ws.Range("T5").Value = 0.4474374422213705
$ws.Range("I6").Value = 0.4564376967244237
$ws.Range("J6").Value = 0.4564376967244237
$ws.Range("O6").Value = 0.009314836229169864
$ws.Range("P6").Value = 0.009314836229169864
$ws.Range("S6").Value = 0.004251642393807508
$ws.Range("T6").Value = 0.004251642393807508
$ws.Range("I7").Value = 0.4564376967244237
$ws.Range("J7").Value = 0.4564376967244237
$ws.Range("O7").Value = 0.01040363699870454
$ws.Range("P7").Value = 0.01040363699870454
$ws.Range("S7").Value = 0.004748612109245697
$ws.Range("T7").Value = 0.004748612109245696
$ws.Range("G8").Value = 0.5724183333333334
$ws.Range("I8").Value = 0.06722654632252778
$ws.Range("J8").Value = 0.06722654632252777
$ws.Range("M8").Value = 14.9419
$ws.Range("N8").Value = 44.8257
$ws.Range("O8").Value = 0.9802815267721257
$ws.Range("P8").Value = 0.9802815267721257
$ws.Range("Q8").Value = 8.553017494833332
$ws.Range("R8").Value = 76.9771574535
$ws.Range("S8").Value = 0.06590094146866456
$ws.Range("T8").Value = 0.06590094146866454
$ws.Range("G9").Value = 0.5724183333333334
$ws.Range("I9").Value = 0.06722654632252778
$ws.Range("J9").Value = 0.06722654632252777
$ws.Range("O9").Value = 0.009314836229169864
$ws.Range("P9").Value = 0.009314836229169864
$ws.Range("Q9").Value = 0.081272527385
$ws.Range("R9").Value = 0.731452746465
$ws.Range("S9").Value = 0.0006262042692470478
$ws.Range("T9").Value = 0.0006262042692470477
$ws.Range("G10").Value = 0.5724183333333334
$ws.Range("I10").Value = 0.06722654632252778
$ws.Range("J10").Value = 0.06722654632252777
$ws.Range("O10").Value = 0.01040363699870454
$ws.Range("P10").Value = 0.01040363699870454
$ws.Range("Q10").Value = 0.090772382045
$ws.Range("R10").Value = 0.8169514384050001
$ws.Range("S10").Value = 0.0006994005846161748
$ws.Range("T10").Value = 0.0006994005846161745